# Automatische test-sync: 2025-06-20 13:00:50
#
# Adds a new incoming-mail log entry (row 12) to the "Logs" sheet and
# updates the "Dashboard" summary sheet / conditional formatting ranges
# accordingly.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Logs sheet: append the new row (row 12)
# ---------------------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A12").Value = "Offerte voor zakelijke samenwerking"
$logs.Range("B12").Value = "mailmind.test@zohomail.eu"
$logs.Range("C12").Value = "Kunt u mij een offerte sturen voor 100 stuks product X?"
$logs.Range("D12").Value = "Offerte / Prijsaanvraag"
$logs.Range("F12").Value = "2025-06-20 13:00:12"
$logs.Range("G12").Value = "Nee"

# ---------------------------------------------------------------------
# 2. Logs sheet: extend the conditional-formatting ranges to row 12
# ---------------------------------------------------------------------
$dFormats = $logs.Range("D2:D11").FormatConditions
$dFormats.Item(1).ModifyAppliesToRange($logs.Range("D2:D12"))

$gFormats = $logs.Range("G2:G11").FormatConditions
$gFormats.Item(1).ModifyAppliesToRange($logs.Range("G2:G12"))

# ---------------------------------------------------------------------
# 3. Dashboard sheet: re-rank the category counts now that
#    "Offerte / Prijsaanvraag" has grown from 1 to 2 occurrences
# ---------------------------------------------------------------------
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Range("A4").Value = "Offerte / Prijsaanvraag"
$dash.Range("B4").Value = 2

$dash.Range("A5").Value = "Productinformatie"
$dash.Range("B5").Value = 1

$dash.Range("A6").Value = "Sollicitatie / Vacature"
$dash.Range("B6").Value = 1
